# Updated cryptos list values (prices / volume(1h)) per latest scrape.
# Values that look like plain decimal numbers (e.g. "1.003") are prefixed
# with a leading apostrophe so Excel stores them as text (matching the
# original inline-string cells) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.896.47'
$ws.Range('D3').Value = '1.633.12'
$ws.Range('E3').Value = '  -2.61%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'209.27"
$ws.Range('E5').Value = '  -1.18%  '
$ws.Range('D6').Value = "'0.5205"
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('D7').Value = "'1.003"
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = "'0.2563"
$ws.Range('E8').Value = '  -3.56%  '
$ws.Range('D9').Value = "'0.06235"
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('D10').Value = "'20.25"
$ws.Range('E10').Value = '  -5.39%  '
$ws.Range('D11').Value = "'0.07562"
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '1.645.48'
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('D13').Value = "'4.354"
$ws.Range('E13').Value = '  -2.20%  '
$ws.Range('D14').Value = '1.859.37'
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').Value = "'0.5412"
$ws.Range('E15').Value = '  -4.13%  '
$ws.Range('D16').Value = '0.0₅7915'
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('D17').Value = "'64.49"
$ws.Range('D18').Value = '25.903.12'
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('D20').Value = "'4.609"
$ws.Range('E20').Value = '  -4.55%  '
$ws.Range('D21').Value = "'184.26"
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('E22').Value = '  -4.20%  '
$ws.Range('D23').Value = "'6.063"
$ws.Range('E23').Value = '  -2.08%  '
$ws.Range('D24').Value = "'1.004"
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = "'145.55"
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('E26').Value = '  -3.98%  '
$ws.Range('D27').Value = "'7.338"
$ws.Range('E27').Value = '  -3.34%  '
$ws.Range('D28').Value = "'15.48"
$ws.Range('E28').Value = '  -3.60%  '
$ws.Range('E29').Value = '  +0.94%  '
$ws.Range('D30').Value = "'0.05943"
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('D32').Value = "'3.345"
$ws.Range('E33').Value = '  -4.49%  '
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('D35').Value = "'0.9691"
$ws.Range('E35').Value = '  -3.36%  '
$ws.Range('D36').Value = "'2.381"
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('D37').Value = "'2.735"
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('D38').Value = "'0.5800"
$ws.Range('E38').Value = '  -4.45%  '
$ws.Range('D39').Value = "'0.01591"
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('D41').Value = "'0.8388"
$ws.Range('E41').Value = '  -3.45%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'5.653"
$ws.Range('E42').Value = '  -7.32%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.020.82'
$ws.Range('E43').Value = '  -5.86%  '
$ws.Range('D44').Value = "'99.54"
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('E46').Value = '  -3.76%  '
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('E48').Value = '  -3.59%  '
$ws.Range('D49').Value = "'7.952"
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('D50').Value = "'0.05173"
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = "'0.4228"
$ws.Range('E51').Value = '  -0.68%  '
